$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points conversion (1 pt = 12700 EMU)
$left   = 141890   / 12700
$top    = 3090041  / 12700
$width  = 5596758  / 12700
$height = 215444   / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 1"

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

$tf.TextFrame.MarginLeft = 0
$tf.TextFrame.MarginTop = 0
$tf.TextFrame.MarginRight = 0
$tf.TextFrame.MarginBottom = 0

$tr = $tf.TextRange
$tr.Text = "date_range"
$tr.Font.Size = 14
$tr.Font.Bold = $false
$tr.Font.Color.ObjectThemeColor = 13
$tr.ParagraphFormat.Alignment = 1

$run1 = $tr.Characters(1, 1)
$run1.Font.Bold = $true
